# Update symbol list (cryptos) - price (D) and 1h volume change (E) columns
# for the rows affected by the Sun Jan 8 03:42:34 UTC 2023 GitHub Actions run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($sheet, $addr, $val) {
    # Force the cell to remain a plain text value (matching the workbook's
    # existing inlineStr/text cells) instead of letting Excel auto-convert
    # numeric- or percentage-looking strings into numbers.
    $cell = $sheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2 (BNB)
Set-TextCell $ws "D2" "260.90"
Set-TextCell $ws "E2" "-0.33%"

# Row 3 (OKB)
Set-TextCell $ws "D3" "27.19"
Set-TextCell $ws "E3" "-0.80%"

# Row 4 (HuobiToken)
Set-TextCell $ws "D4" "4.710"
Set-TextCell $ws "E4" "-1.22%"

# Row 5 (Cronos)
Set-TextCell $ws "D5" "0.06227"
Set-TextCell $ws "E5" "2.57%"

# Row 6 (KuCoinToken)
Set-TextCell $ws "D6" "6.731"
Set-TextCell $ws "E6" "0.26%"

# Row 7 (MXToken)
Set-TextCell $ws "D7" "0.8517"
Set-TextCell $ws "E7" "-1.32%"

# Row 8 (FTXToken)
Set-TextCell $ws "D8" "0.9095"
Set-TextCell $ws "E8" "-1.52%"

# Row 9 (WazirX) - only E changes
Set-TextCell $ws "E9" "-0.86%"

# Row 10 (LiechtensteinCryptoassetsExchange)
Set-TextCell $ws "D10" "0.04819"
Set-TextCell $ws "E10" "-3.92%"

# Row 11 (MandalaExchangeToken)
Set-TextCell $ws "D11" "0.07089"
Set-TextCell $ws "E11" "-0.96%"

# Row 12 (BitrueCoin) - only E changes
Set-TextCell $ws "E12" "2.26%"

# Row 13 (BitMartToken)
Set-TextCell $ws "D13" "0.09053"
Set-TextCell $ws "E13" "-0.63%"

# Row 14 (BitForexToken)
Set-TextCell $ws "D14" "0.001537"
Set-TextCell $ws "E14" "-0.38%"

# Row 15 (One)
Set-TextCell $ws "D15" "0.0006166"
Set-TextCell $ws "E15" "1.38%"

# Row 16 (TigerCash)
Set-TextCell $ws "D16" "0.006064"
Set-TextCell $ws "E16" "-2.06%"

# Row 17 (LEO)
Set-TextCell $ws "D17" "3.453"
Set-TextCell $ws "E17" "0.06%"

# Row 18 (GateToken)
Set-TextCell $ws "D18" "3.173"
Set-TextCell $ws "E18" "0.16%"

# Row 19 (BTSEToken) - only E changes
Set-TextCell $ws "E19" "-0.55%"

# Row 21 (ProBitToken) - only E changes
Set-TextCell $ws "E21" "0.91%"

# Row 22 (MCDex)
Set-TextCell $ws "D22" "4.101"
Set-TextCell $ws "E22" "0.06%"

# Row 23 (CoinExToken)
Set-TextCell $ws "D23" "0.04232"
Set-TextCell $ws "E23" "-0.79%"

# Row 24 (BitKan)
Set-TextCell $ws "D24" "0.001221"
Set-TextCell $ws "E24" "0.34%"

# Row 26 (NitroEx)
Set-TextCell $ws "D26" "0.0001200"
Set-TextCell $ws "E26" "0.00%"

# Row 40 (IDEX)
Set-TextCell $ws "D40" "0.03919"
Set-TextCell $ws "E40" "0.95%"

# Row 41 (BKEXToken)
Set-TextCell $ws "D41" "0.1112"
Set-TextCell $ws "E41" "-0.29%"

# Row 42 (KickToken) - only D changes
Set-TextCell $ws "D42" "0.004113"

# Row 43 (CEJI)
Set-TextCell $ws "D43" "0.002152"
Set-TextCell $ws "E43" "-2.56%"

# Row 44 (LocalTraders)
Set-TextCell $ws "D44" "0.01388"
Set-TextCell $ws "E44" "-7.16%"

# Row 45 (CoinLion)
Set-TextCell $ws "D45" "0.00005105"
Set-TextCell $ws "E45" "-3.64%"

# Row 47 (CoinbaseStockToken) - only D changes
Set-TextCell $ws "D47" "0.03402"

# Row 48 (BOLO)
Set-TextCell $ws "D48" "0.06499"
Set-TextCell $ws "E48" "-50.81%"
